$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing history at row 254,
# pushing every subsequent row (old 254..355) down by one (new 255..356).
$ws.Rows.Item(254).Insert()

$ws.Cells.Item(254,1).Value  = 4
$ws.Cells.Item(254,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(254,3).Value  = "Los Lagos"
$ws.Cells.Item(254,4).Value  = 44875
$ws.Cells.Item(254,5).Value  = 10
$ws.Cells.Item(254,6).Value  = 100112037
$ws.Cells.Item(254,7).Value  = "Cebollín"
$ws.Cells.Item(254,8).Value  = "Sin especificar"
$ws.Cells.Item(254,9).Value  = "Primera"
$ws.Cells.Item(254,10).Value = 120
$ws.Cells.Item(254,11).Value = 6500
$ws.Cells.Item(254,12).Value = 6500
$ws.Cells.Item(254,13).Value = 6500
$ws.Cells.Item(254,14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(254,15).Value = "Región Metropolitana"
$ws.Cells.Item(254,16).Value = 181
$ws.Cells.Item(254,17).Value = 36
$ws.Cells.Item(254,18).Value = "Hortaliza"

# Keep the same date display/format (s="2" / numFmt 165) as the other date cells in column D
$ws.Cells.Item(254,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
